{"js": "// Update the division-problem worksheet table: replace each problem's\n// text with the new problem, cell-by-cell (by row/column position) so\n// that duplicate \"before\" values (e.g. \"37\u00f72=\" appears twice) are each\n// mapped to their own distinct \"after\" value.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index -> [oldText -> newText] in column order (cols 0..4)\nconst rowUpdates = {\n  0: [\"66\u00f77=\", \"51\u00f76=\", \"21\u00f78=\", \"47\u00f73=\", \"95\u00f77=\"],\n  4: [\"55\u00f73=\", \"58\u00f77=\", \"47\u00f79=\", \"96\u00f75=\", \"92\u00f79=\"],\n  8: [\"66\u00f73=\", \"47\u00f74=\", \"41\u00f74=\", \"57\u00f74=\", \"51\u00f74=\"],\n  12: [\"25\u00f72=\", \"60\u00f74=\", \"70\u00f77=\", \"30\u00f72=\", \"86\u00f73=\"],\n  16: [\"15\u00f79=\", \"78\u00f76=\", \"53\u00f73=\", \"88\u00f79=\", \"60\u00f76=\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const newValues = rowUpdates[rowIndex];\n  for (let col = 0; col < newValues.length; col++) {\n    const cell = table.getCell(parseInt(rowIndex, 10), col);\n    cell.value = newValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: replace each problem's\n# text with the new problem, cell-by-cell (by row/column position) so\n# that duplicate \"before\" values (e.g. \"37\u00f72=\" appears twice) are each\n# mapped to their own distinct \"after\" value.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"66\u00f77=\", \"51\u00f76=\", \"21\u00f78=\", \"47\u00f73=\", \"95\u00f77=\")\n    5  = @(\"55\u00f73=\", \"58\u00f77=\", \"47\u00f79=\", \"96\u00f75=\", \"92\u00f79=\")\n    9  = @(\"66\u00f73=\", \"47\u00f74=\", \"41\u00f74=\", \"57\u00f74=\", \"51\u00f74=\")\n    13 = @(\"25\u00f72=\", \"60\u00f74=\", \"70\u00f77=\", \"30\u00f72=\", \"86\u00f73=\")\n    17 = @(\"15\u00f79=\", \"78\u00f76=\", \"53\u00f73=\", \"88\u00f79=\", \"60\u00f76=\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $newValues = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $newValues.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $newValues[$col - 1]\n    }\n}\n"}
